$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.735.66"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "3.454.88"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'573.95"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "'159.38"
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.453.33"
$ws.Range("D9").Value = "'0.581"
$ws.Range("E9").Value = "  -5.42%  "
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("E11").Value = "  -1.83%  "
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("D13").Value = "4.044.91"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").Value = "'27.54"
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").Value = "'0.0000175"
$ws.Range("E16").Value = "  -9.54%  "
$ws.Range("D17").Value = "64.743.48"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "3.457.60"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("E19").Value = "  -3.85%  "
$ws.Range("D20").Value = "'13.81"
$ws.Range("E20").Value = "  -3.68%  "
$ws.Range("D21").Value = "'380.73"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").Value = "'8.00"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "'72.40"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  -3.13%  "
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("E27").Value = "  -1.57%  "
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("E30").Value = "  -3.98%  "
$ws.Range("D31").Value = "'6.08"
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("E32").Value = "  -1.30%  "
$ws.Range("D33").Value = "'23.28"
$ws.Range("E33").Value = "  -1.34%  "
$ws.Range("E34").Value = "  -1.20%  "
$ws.Range("E35").Value = "  -2.23%  "
$ws.Range("D36").Value = "'161.53"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("E37").Value = "  -2.08%  "
$ws.Range("D38").Value = "2.880.99"
$ws.Range("D39").Value = "'0.0749"
$ws.Range("E39").Value = "  -3.56%  "
$ws.Range("D40").Value = "'26.55"
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("D41").Value = "'0.802"
$ws.Range("E41").Value = "  +3.59%  "
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("D43").Value = "'43.04"
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("D44").Value = "'6.52"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("E45").Value = "  +2.22%  "
$ws.Range("E46").Value = "  -2.42%  "
$ws.Range("E47").Value = "  +12.67%  "
$ws.Range("D48").Value = "'323.40"
$ws.Range("E48").Value = "  +4.17%  "
$ws.Range("E49").Value = "  -1.85%  "
$ws.Range("E50").Value = "  -2.44%  "
$ws.Range("E51").Value = "  -2.10%  "
